$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as at" date in the intro paragraph (A2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 14 March 2025"

# Drop the oldest week (row 5) by shifting every subsequent row up by one,
# then removing the now-duplicated last row (69).
# Row 5 <- old row 6
$ws.Range("A5").Value = "17 Mar 2025"
$ws.Range("B5").Value = "A Profile of Repeat Offending by Children and Young People in England and Wales"
$ws.Range("C5").Value = "20 March 2025"
$ws.Range("D5").Value = "confirmed"
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = "standard"

# Row 6 <- old row 7
$ws.Range("A6").Value = "24 Mar 2025"
$ws.Range("B6").Value = "Criminal court statistics quarterly: October to December 2024"
$ws.Range("C6").Value = "27 March 2025"
$ws.Range("D6").Value = "provisional"
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = "standard"

# Row 7 <- old row 8
$ws.Range("A7").Value = "24 Mar 2025"
$ws.Range("B7").Value = "Family court statistics quarterly: October to December 2024 "
$ws.Range("C7").Value = "27 March 2025"
$ws.Range("D7").Value = "provisional"
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = "standard"

# Row 8 <- old row 9
$ws.Range("A8").Value = "24 Mar 2025"
$ws.Range("B8").Value = "Legal aid statistics quarterly: October to December 2024"
$ws.Range("C8").Value = "27 March 2025"
$ws.Range("D8").Value = "confirmed"
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = "standard"

# Row 9 <- old row 10
$ws.Range("A9").Value = "31 Mar 2025"
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = $null
$ws.Range("E9").Value = 14
$ws.Range("F9").Value = $null

# Row 10 <- old row 11
$ws.Range("A10").Value = "07 Apr 2025"
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = $null

# Row 11 <- old row 12
$ws.Range("A11").Value = "14 Apr 2025"
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = $null

# Row 12 <- old row 13
$ws.Range("A12").Value = "21 Apr 2025"
$ws.Range("B12").Value = "Safety in the children and young people secure estate: Update to December 2024"
$ws.Range("C12").Value = "24 April 2025"
$ws.Range("D12").Value = "confirmed"
$ws.Range("E12").Value = 17
$ws.Range("F12").Value = "standard"

# Row 13 <- old row 14
$ws.Range("A13").Value = "21 Apr 2025"
$ws.Range("B13").Value = "Safety in custody: quarterly update to December 2024"
$ws.Range("C13").Value = "24 April 2025"
$ws.Range("D13").Value = "confirmed"
$ws.Range("E13").Value = 17
$ws.Range("F13").Value = "standard"

# Row 14 <- old row 15
$ws.Range("A14").Value = "21 Apr 2025"
$ws.Range("B14").Value = "Justice data lab statistics: April 2025"
$ws.Range("C14").Value = "24 April 2025"
$ws.Range("D14").Value = "provisional"
$ws.Range("E14").Value = 17
$ws.Range("F14").Value = "standard"

# Row 15 <- old row 16
$ws.Range("A15").Value = "21 Apr 2025"
$ws.Range("B15").Value = "Proven reoffending statistics: April to June 2023"
$ws.Range("C15").Value = "24 April 2025"
$ws.Range("D15").Value = "provisional"
$ws.Range("E15").Value = 17
$ws.Range("F15").Value = "standard"

# Row 16 <- old row 17
$ws.Range("A16").Value = "21 Apr 2025"
$ws.Range("B16").Value = "Safety in the children and young people secure estate: Update to December 2024"
$ws.Range("C16").Value = "24 April 2025"
$ws.Range("D16").Value = "provisional"
$ws.Range("E16").Value = 17
$ws.Range("F16").Value = "standard"

# Row 17 <- old row 18
$ws.Range("A17").Value = "21 Apr 2025"
$ws.Range("B17").Value = "Offender management statistics quarterly: October to December 2024"
$ws.Range("C17").Value = "24 April 2025"
$ws.Range("D17").Value = "provisional"
$ws.Range("E17").Value = 17
$ws.Range("F17").Value = "standard"

# Row 18 <- old row 19
$ws.Range("A18").Value = "28 Apr 2025"
$ws.Range("B18").Value = $null
$ws.Range("C18").Value = $null
$ws.Range("D18").Value = $null
$ws.Range("E18").Value = 18
$ws.Range("F18").Value = $null

# Row 19 <- old row 20
$ws.Range("A19").Value = "05 May 2025"
$ws.Range("B19").Value = "Coroners statistics 2024"
$ws.Range("C19").Value = "8 May 2025"
$ws.Range("D19").Value = "provisional"
$ws.Range("E19").Value = 19
$ws.Range("F19").Value = "standard"

# Row 20 <- old row 21
$ws.Range("A20").Value = "12 May 2025"
$ws.Range("B20").Value = "HM Prison and Probation Service workforce quarterly: March 2025."
$ws.Range("C20").Value = "15 May 2025"
$ws.Range("D20").Value = "confirmed"
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = "standard"

# Row 21 <- old row 22
$ws.Range("A21").Value = "12 May 2025"
$ws.Range("B21").Value = "Mortgage and landlord possession statistics: January to March 2025"
$ws.Range("C21").Value = "15 May 2025"
$ws.Range("D21").Value = "provisional"
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = "standard"

# Row 22 <- old row 23
$ws.Range("A22").Value = "12 May 2025"
$ws.Range("B22").Value = "First time entrants (FTE) into the Criminal Justice System and Offender Histories: year ending December 2024"
$ws.Range("C22").Value = "15 May 2025"
$ws.Range("D22").Value = "provisional"
$ws.Range("E22").Value = 20
$ws.Range("F22").Value = "standard"

# Row 23 <- old row 24
$ws.Range("A23").Value = "12 May 2025"
$ws.Range("B23").Value = "Knife and Offensive Weapon Sentencing Statistics:  October to December 2024"
$ws.Range("C23").Value = "15 May 2025"
$ws.Range("D23").Value = "provisional"
$ws.Range("E23").Value = 20
$ws.Range("F23").Value = "standard"

# Row 24 <- old row 25
$ws.Range("A24").Value = "12 May 2025"
$ws.Range("B24").Value = "First time entrants (FTE) into the Criminal Justice System and Offender Histories: year ending December 2024"
$ws.Range("C24").Value = "15 May 2025"
$ws.Range("D24").Value = "confirmed"
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = "standard"

# Row 25 <- old row 26
$ws.Range("A25").Value = "12 May 2025"
$ws.Range("B25").Value = "Knife and Offensive Weapon Sentencing Statistics:  October to December 2024"
$ws.Range("C25").Value = "15 May 2025"
$ws.Range("D25").Value = "confirmed"
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = "standard"

# Row 26 <- old row 27
$ws.Range("A26").Value = "19 May 2025"
$ws.Range("B26").Value = $null
$ws.Range("C26").Value = $null
$ws.Range("D26").Value = $null
$ws.Range("E26").Value = 21
$ws.Range("F26").Value = $null

# Row 27 <- old row 28
$ws.Range("A27").Value = "26 May 2025"
$ws.Range("B27").Value = $null
$ws.Range("C27").Value = $null
$ws.Range("D27").Value = $null
$ws.Range("E27").Value = 22
$ws.Range("F27").Value = $null

# Row 28 <- old row 29
$ws.Range("A28").Value = "02 Jun 2025"
$ws.Range("B28").Value = "Civil justice statistics: January to March 2025"
$ws.Range("C28").Value = "5 June 2025"
$ws.Range("D28").Value = "confirmed"
$ws.Range("E28").Value = 23
$ws.Range("F28").Value = "standard"

# Row 29 <- old row 30
$ws.Range("A29").Value = "09 Jun 2025"
$ws.Range("B29").Value = "Tribunals statistics quarterly: January to March 2025"
$ws.Range("C29").Value = "12 June 2025"
$ws.Range("D29").Value = "confirmed"
$ws.Range("E29").Value = 24
$ws.Range("F29").Value = "standard"

# Row 30 <- old row 31
$ws.Range("A30").Value = "16 Jun 2025"
$ws.Range("B30").Value = $null
$ws.Range("C30").Value = $null
$ws.Range("D30").Value = $null
$ws.Range("E30").Value = 25
$ws.Range("F30").Value = $null

# Row 31 <- old row 32
$ws.Range("A31").Value = "23 Jun 2025"
$ws.Range("B31").Value = "Legal aid statistics quarterly: January to March 2025"
$ws.Range("C31").Value = "26 June 2025"
$ws.Range("D31").Value = "provisional"
$ws.Range("E31").Value = 26
$ws.Range("F31").Value = "standard"

# Row 32 <- old row 33
$ws.Range("A32").Value = "23 Jun 2025"
$ws.Range("B32").Value = "Family court statistics quarterly: January to March 2025"
$ws.Range("C32").Value = "26 June 2025"
$ws.Range("D32").Value = "provisional"
$ws.Range("E32").Value = 26
$ws.Range("F32").Value = "standard"

# Row 33 <- old row 34
$ws.Range("A33").Value = "30 Jun 2025"
$ws.Range("B33").Value = $null
$ws.Range("C33").Value = $null
$ws.Range("D33").Value = $null
$ws.Range("E33").Value = 27
$ws.Range("F33").Value = $null

# Row 34 <- old row 35
$ws.Range("A34").Value = "07 Jul 2025"
$ws.Range("B34").Value = "Diversity of the judiciary: 2025 statistics"
$ws.Range("C34").Value = "10 July 2025"
$ws.Range("D34").Value = "provisional"
$ws.Range("E34").Value = 28
$ws.Range("F34").Value = "standard"

# Row 35 <- old row 36
$ws.Range("A35").Value = "14 Jul 2025"
$ws.Range("B35").Value = $null
$ws.Range("C35").Value = $null
$ws.Range("D35").Value = $null
$ws.Range("E35").Value = 29
$ws.Range("F35").Value = $null

# Row 36 <- old row 37
$ws.Range("A36").Value = "21 Jul 2025"
$ws.Range("B36").Value = $null
$ws.Range("C36").Value = $null
$ws.Range("D36").Value = $null
$ws.Range("E36").Value = 30
$ws.Range("F36").Value = $null

# Row 37 <- old row 38
$ws.Range("A37").Value = "28 Jul 2025"
$ws.Range("B37").Value = "Safety in custody: quarterly update to March 2025"
$ws.Range("C37").Value = "31 July 2025"
$ws.Range("D37").Value = "provisional"
$ws.Range("E37").Value = 31
$ws.Range("F37").Value = "standard"

# Row 38 <- old row 39
$ws.Range("A38").Value = "28 Jul 2025"
$ws.Range("B38").Value = "Offender management statistics quarterly: January 2025 to March 2025 "
$ws.Range("C38").Value = "31 July 2025"
$ws.Range("D38").Value = "provisional"
$ws.Range("E38").Value = 31
$ws.Range("F38").Value = "standard"

# Row 39 <- old row 40
$ws.Range("A39").Value = "28 Jul 2025"
$ws.Range("B39").Value = "Proven reoffending statistics: July 2023 and September 2023"
$ws.Range("C39").Value = "31 July 2025"
$ws.Range("D39").Value = "provisional"
$ws.Range("E39").Value = 31
$ws.Range("F39").Value = "standard"

# Row 40 <- old row 41
$ws.Range("A40").Value = "28 Jul 2025"
$ws.Range("B40").Value = "Safety in the children and young people secure estate: Update to March 2025"
$ws.Range("C40").Value = "31 July 2025"
$ws.Range("D40").Value = "provisional"
$ws.Range("E40").Value = 31
$ws.Range("F40").Value = "standard"

# Row 41 <- old row 42
$ws.Range("A41").Value = "28 Jul 2025"
$ws.Range("B41").Value = "Prison Performance Ratings: 2024 to 2025"
$ws.Range("C41").Value = "31 July 2025"
$ws.Range("D41").Value = "provisional"
$ws.Range("E41").Value = 31
$ws.Range("F41").Value = "standard"

# Row 42 <- old row 43
$ws.Range("A42").Value = "04 Aug 2025"
$ws.Range("B42").Value = $null
$ws.Range("C42").Value = $null
$ws.Range("D42").Value = $null
$ws.Range("E42").Value = 32
$ws.Range("F42").Value = $null

# Row 43 <- old row 44
$ws.Range("A43").Value = "11 Aug 2025"
$ws.Range("B43").Value = "Mortgage and landlord possession statistics:  April to June 2025"
$ws.Range("C43").Value = "14 August 2025"
$ws.Range("D43").Value = "provisional"
$ws.Range("E43").Value = 33
$ws.Range("F43").Value = "standard"

# Row 44 <- old row 45
$ws.Range("A44").Value = "18 Aug 2025"
$ws.Range("B44").Value = "Knife and Offensive Weapon Sentencing Statistics:  January to March 2025"
$ws.Range("C44").Value = "21 August 2025"
$ws.Range("D44").Value = "provisional"
$ws.Range("E44").Value = 34
$ws.Range("F44").Value = "standard"

# Row 45 <- old row 46
$ws.Range("A45").Value = "18 Aug 2025"
$ws.Range("B45").Value = "HM Prison and Probation Service workforce quarterly: June 2025"
$ws.Range("C45").Value = "21 August 2025"
$ws.Range("D45").Value = "provisional"
$ws.Range("E45").Value = 34
$ws.Range("F45").Value = "standard"

# Row 46 <- old row 47
$ws.Range("A46").Value = "25 Aug 2025"
$ws.Range("B46").Value = $null
$ws.Range("C46").Value = $null
$ws.Range("D46").Value = $null
$ws.Range("E46").Value = 35
$ws.Range("F46").Value = $null

# Row 47 <- old row 48
$ws.Range("A47").Value = "01 Sep 2025"
$ws.Range("B47").Value = "Civil justice statistics: April to June 2025"
$ws.Range("C47").Value = "4 September 2025"
$ws.Range("D47").Value = "provisional"
$ws.Range("E47").Value = 36
$ws.Range("F47").Value = "standard"

# Row 48 <- old row 49
$ws.Range("A48").Value = "08 Sep 2025"
$ws.Range("B48").Value = "Tribunals statistics quarterly: April to June 2025 "
$ws.Range("C48").Value = "11 September 2025"
$ws.Range("D48").Value = "provisional"
$ws.Range("E48").Value = 37
$ws.Range("F48").Value = "standard"

# Row 49 <- old row 50
$ws.Range("A49").Value = "15 Sep 2025"
$ws.Range("B49").Value = $null
$ws.Range("C49").Value = $null
$ws.Range("D49").Value = $null
$ws.Range("E49").Value = 38
$ws.Range("F49").Value = $null

# Row 50 <- old row 51
$ws.Range("A50").Value = "22 Sep 2025"
$ws.Range("B50").Value = "Family court statistics quarterly: April to June 2025"
$ws.Range("C50").Value = "25 September 2025"
$ws.Range("D50").Value = "provisional"
$ws.Range("E50").Value = 39
$ws.Range("F50").Value = "standard"

# Row 51 <- old row 52
$ws.Range("A51").Value = "22 Sep 2025"
$ws.Range("B51").Value = "Legal aid statistics quarterly: April to June 2025"
$ws.Range("C51").Value = "25 September 2025"
$ws.Range("D51").Value = "provisional"
$ws.Range("E51").Value = 39
$ws.Range("F51").Value = "standard"

# Row 52 <- old row 53
$ws.Range("A52").Value = "29 Sep 2025"
$ws.Range("B52").Value = $null
$ws.Range("C52").Value = $null
$ws.Range("D52").Value = $null
$ws.Range("E52").Value = 40
$ws.Range("F52").Value = $null

# Row 53 <- old row 54
$ws.Range("A53").Value = "06 Oct 2025"
$ws.Range("B53").Value = $null
$ws.Range("C53").Value = $null
$ws.Range("D53").Value = $null
$ws.Range("E53").Value = 41
$ws.Range("F53").Value = $null

# Row 54 <- old row 55
$ws.Range("A54").Value = "13 Oct 2025"
$ws.Range("B54").Value = $null
$ws.Range("C54").Value = $null
$ws.Range("D54").Value = $null
$ws.Range("E54").Value = 42
$ws.Range("F54").Value = $null

# Row 55 <- old row 56
$ws.Range("A55").Value = "20 Oct 2025"
$ws.Range("B55").Value = $null
$ws.Range("C55").Value = $null
$ws.Range("D55").Value = $null
$ws.Range("E55").Value = 43
$ws.Range("F55").Value = $null

# Row 56 <- old row 57
$ws.Range("A56").Value = "27 Oct 2025"
$ws.Range("B56").Value = "Safety in custody: quarterly update to June 2025"
$ws.Range("C56").Value = "30 October 2025"
$ws.Range("D56").Value = "provisional"
$ws.Range("E56").Value = 44
$ws.Range("F56").Value = "standard"

# Row 57 <- old row 58
$ws.Range("A57").Value = "27 Oct 2025"
$ws.Range("B57").Value = "Proven reoffending statistics: October to December 2023"
$ws.Range("C57").Value = "30 October 2025"
$ws.Range("D57").Value = "provisional"
$ws.Range("E57").Value = 44
$ws.Range("F57").Value = "standard"

# Row 58 <- old row 59
$ws.Range("A58").Value = "27 Oct 2025"
$ws.Range("B58").Value = "Deaths of offenders supervised in the community, England and Wales, 2024/2025"
$ws.Range("C58").Value = "30 October 2025"
$ws.Range("D58").Value = "provisional"
$ws.Range("E58").Value = 44
$ws.Range("F58").Value = "standard"

# Row 59 <- old row 60
$ws.Range("A59").Value = "03 Nov 2025"
$ws.Range("B59").Value = $null
$ws.Range("C59").Value = $null
$ws.Range("D59").Value = $null
$ws.Range("E59").Value = 45
$ws.Range("F59").Value = $null

# Row 60 <- old row 61
$ws.Range("A60").Value = "10 Nov 2025"
$ws.Range("B60").Value = $null
$ws.Range("C60").Value = $null
$ws.Range("D60").Value = $null
$ws.Range("E60").Value = 46
$ws.Range("F60").Value = $null

# Row 61 <- old row 62
$ws.Range("A61").Value = "17 Nov 2025"
$ws.Range("B61").Value = "Knife and Offensive Weapon Sentencing Statistics:  April to June 2025"
$ws.Range("C61").Value = "20 November 2025"
$ws.Range("D61").Value = "provisional"
$ws.Range("E61").Value = 47
$ws.Range("F61").Value = "standard"

# Row 62 <- old row 63
$ws.Range("A62").Value = "17 Nov 2025"
$ws.Range("B62").Value = " HM Prison and Probation Service workforce quarterly: September 2025"
$ws.Range("C62").Value = "20 November 2025"
$ws.Range("D62").Value = "provisional"
$ws.Range("E62").Value = 47
$ws.Range("F62").Value = "standard"

# Row 63 <- old row 64
$ws.Range("A63").Value = "24 Nov 2025"
$ws.Range("B63").Value = "Her Majesty’s Prison and Probation Service offender equalities report: 2024 to 2025"
$ws.Range("C63").Value = "27 November 2025"
$ws.Range("D63").Value = "provisional"
$ws.Range("E63").Value = 48
$ws.Range("F63").Value = "standard"

# Row 64 <- old row 65
$ws.Range("A64").Value = "01 Dec 2025"
$ws.Range("B64").Value = " Civil justice statistics: July to September 2025"
$ws.Range("C64").Value = "4 December 2025"
$ws.Range("D64").Value = "provisional"
$ws.Range("E64").Value = 49
$ws.Range("F64").Value = "standard"

# Row 65 <- old row 66
$ws.Range("A65").Value = "08 Dec 2025"
$ws.Range("B65").Value = "Tribunals statistics quarterly: July to September 2025"
$ws.Range("C65").Value = "11 December 2025"
$ws.Range("D65").Value = "provisional"
$ws.Range("E65").Value = 50
$ws.Range("F65").Value = "standard"

# Row 66 <- old row 67
$ws.Range("A66").Value = "15 Dec 2025"
$ws.Range("B66").Value = "Family court statistics quarterly: July to September 2025"
$ws.Range("C66").Value = "18 December 2025"
$ws.Range("D66").Value = "provisional"
$ws.Range("E66").Value = 51
$ws.Range("F66").Value = "standard"

# Row 67 <- old row 68
$ws.Range("A67").Value = "15 Dec 2025"
$ws.Range("B67").Value = "Criminal court statistics quarterly: July to September 2025"
$ws.Range("C67").Value = "18 December 2025"
$ws.Range("D67").Value = "provisional"
$ws.Range("E67").Value = 51
$ws.Range("F67").Value = "standard"

# Row 68 <- old row 69
$ws.Range("A68").Value = "15 Dec 2025"
$ws.Range("B68").Value = "Legal aid statistics quarterly: July to September 2025"
$ws.Range("C68").Value = "18 December 2025"
$ws.Range("D68").Value = "provisional"
$ws.Range("E68").Value = 51
$ws.Range("F68").Value = "standard"

# Remove the now-empty trailing row
$ws.Range("A69:F69").EntireRow.Delete()

# Shrink the conditional formatting ranges to match the reduced data range
$fcs = $ws.Range("A5").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$A`$5:`$F`$69") {
        $fc.ModifyAppliesToRange($ws.Range("A5:F68"))
    } elseif ($addr -eq "`$A`$5:`$A`$69") {
        $fc.ModifyAppliesToRange($ws.Range("A5:A68"))
    }
}

Write-Host "done"